# Skill.xlsx -- unify the conception of DataNode, DataTable, Entity.
#
# The sheet that used to represent a generic "Property" table is renamed
# to "DataNode" to match the new naming convention, and the active
# selection on the sheet's frozen (bottom-left) pane is moved from K9 to
# O40 (reflecting where the author was last working in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Property1" -> "DataNode"
$ws.Name = "DataNode"

# Update the active cell/selection in the frozen pane from K9 to O40
$ws.Range("O40").Select()
